$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether it is the numeric-looking Price column (D)
$updates = @(
    @{ Cell = "D2"; Value = "28.203.50"; IsPrice = $true },
    @{ Cell = "E2"; Value = "  +0.19%  "; IsPrice = $false },
    @{ Cell = "D3"; Value = "1.854.54"; IsPrice = $true },
    @{ Cell = "E3"; Value = "  -0.64%  "; IsPrice = $false },
    @{ Cell = "D4"; Value = "1.002"; IsPrice = $true },
    @{ Cell = "E4"; Value = "  +0.22%  "; IsPrice = $false },
    @{ Cell = "D5"; Value = "329.64"; IsPrice = $true },
    @{ Cell = "E5"; Value = "  -2.13%  "; IsPrice = $false },
    @{ Cell = "D6"; Value = "1.002"; IsPrice = $true },
    @{ Cell = "E6"; Value = "  +0.22%  "; IsPrice = $false },
    @{ Cell = "D7"; Value = "0.4538"; IsPrice = $true },
    @{ Cell = "E7"; Value = "  -3.54%  "; IsPrice = $false },
    @{ Cell = "D8"; Value = "0.3915"; IsPrice = $true },
    @{ Cell = "E8"; Value = "  -0.02%  "; IsPrice = $false },
    @{ Cell = "D9"; Value = "47.36"; IsPrice = $true },
    @{ Cell = "E9"; Value = "  +1.04%  "; IsPrice = $false },
    @{ Cell = "D10"; Value = "0.07776"; IsPrice = $true },
    @{ Cell = "E10"; Value = "  -2.51%  "; IsPrice = $false },
    @{ Cell = "D11"; Value = "0.9838"; IsPrice = $true },
    @{ Cell = "E11"; Value = "  -1.26%  "; IsPrice = $false },
    @{ Cell = "D12"; Value = "21.19"; IsPrice = $true },
    @{ Cell = "E12"; Value = "  -2.26%  "; IsPrice = $false },
    @{ Cell = "D13"; Value = "1.857.05"; IsPrice = $true },
    @{ Cell = "E13"; Value = "  +0.42%  "; IsPrice = $false },
    @{ Cell = "D14"; Value = "5.776"; IsPrice = $true },
    @{ Cell = "E14"; Value = "  -3.39%  "; IsPrice = $false },
    @{ Cell = "D15"; Value = "6.933"; IsPrice = $true },
    @{ Cell = "E15"; Value = "  -4.60%  "; IsPrice = $false },
    @{ Cell = "D16"; Value = "1.002"; IsPrice = $true },
    @{ Cell = "E16"; Value = "  +0.00%  "; IsPrice = $false },
    @{ Cell = "D17"; Value = "87.55"; IsPrice = $true },
    @{ Cell = "E17"; Value = "  -4.20%  "; IsPrice = $false },
    @{ Cell = "D18"; Value = "0.06518"; IsPrice = $true },
    @{ Cell = "E18"; Value = "  -1.27%  "; IsPrice = $false },
    @{ Cell = "D19"; Value = "0.00001013"; IsPrice = $true },
    @{ Cell = "E19"; Value = "  -2.95%  "; IsPrice = $false },
    @{ Cell = "D20"; Value = "16.99"; IsPrice = $true },
    @{ Cell = "E20"; Value = "  -4.25%  "; IsPrice = $false },
    @{ Cell = "D21"; Value = "1.001"; IsPrice = $true },
    @{ Cell = "E21"; Value = "  +0.19%  "; IsPrice = $false },
    @{ Cell = "D22"; Value = "28.224.17"; IsPrice = $true },
    @{ Cell = "E22"; Value = "  +0.27%  "; IsPrice = $false },
    @{ Cell = "D23"; Value = "5.271"; IsPrice = $true },
    @{ Cell = "E23"; Value = "  -2.99%  "; IsPrice = $false },
    @{ Cell = "D24"; Value = "10.63"; IsPrice = $true },
    @{ Cell = "E24"; Value = "  -3.58%  "; IsPrice = $false },
    @{ Cell = "D25"; Value = "2.250"; IsPrice = $true },
    @{ Cell = "E25"; Value = "  -1.50%  "; IsPrice = $false },
    @{ Cell = "D26"; Value = "2.074.84"; IsPrice = $true },
    @{ Cell = "E26"; Value = "  +0.32%  "; IsPrice = $false },
    @{ Cell = "D27"; Value = "156.85"; IsPrice = $true },
    @{ Cell = "E27"; Value = "  -1.39%  "; IsPrice = $false },
    @{ Cell = "E28"; Value = "  -3.42%  "; IsPrice = $false },
    @{ Cell = "D29"; Value = "2.036"; IsPrice = $true },
    @{ Cell = "E29"; Value = "  -4.24%  "; IsPrice = $false },
    @{ Cell = "D30"; Value = "5.256"; IsPrice = $true },
    @{ Cell = "E30"; Value = "  -4.41%  "; IsPrice = $false },
    @{ Cell = "D31"; Value = "116.04"; IsPrice = $true },
    @{ Cell = "E31"; Value = "  -3.06%  "; IsPrice = $false },
    @{ Cell = "B32"; Value = "ImmutableX"; IsPrice = $false },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; IsPrice = $false },
    @{ Cell = "D32"; Value = "0.9348"; IsPrice = $true },
    @{ Cell = "E32"; Value = "  -3.96%  "; IsPrice = $false },
    @{ Cell = "B33"; Value = "Stellar"; IsPrice = $false },
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; IsPrice = $false },
    @{ Cell = "D33"; Value = "0.09245"; IsPrice = $true },
    @{ Cell = "E33"; Value = "  -2.58%  "; IsPrice = $false },
    @{ Cell = "D34"; Value = "3.603"; IsPrice = $true },
    @{ Cell = "E34"; Value = "  +0.92%  "; IsPrice = $false },
    @{ Cell = "D35"; Value = "1.371"; IsPrice = $true },
    @{ Cell = "E35"; Value = "  -0.15%  "; IsPrice = $false },
    @{ Cell = "E36"; Value = "  -3.12%  "; IsPrice = $false },
    @{ Cell = "D37"; Value = "0.05984"; IsPrice = $true },
    @{ Cell = "E37"; Value = "  -2.00%  "; IsPrice = $false },
    @{ Cell = "D38"; Value = "0.02188"; IsPrice = $true },
    @{ Cell = "E38"; Value = "  -3.63%  "; IsPrice = $false },
    @{ Cell = "D39"; Value = "8.161"; IsPrice = $true },
    @{ Cell = "E39"; Value = "  -2.88%  "; IsPrice = $false },
    @{ Cell = "D40"; Value = "1.159"; IsPrice = $true },
    @{ Cell = "E40"; Value = "  -1.27%  "; IsPrice = $false },
    @{ Cell = "D41"; Value = "1.001"; IsPrice = $true },
    @{ Cell = "E41"; Value = "  +0.12%  "; IsPrice = $false },
    @{ Cell = "D42"; Value = "0.5641"; IsPrice = $true },
    @{ Cell = "E42"; Value = "  -5.56%  "; IsPrice = $false },
    @{ Cell = "D43"; Value = "9.921"; IsPrice = $true },
    @{ Cell = "E43"; Value = "  -3.97%  "; IsPrice = $false },
    @{ Cell = "D44"; Value = "0.1783"; IsPrice = $true },
    @{ Cell = "E44"; Value = "  -5.21%  "; IsPrice = $false },
    @{ Cell = "D45"; Value = "1.258"; IsPrice = $true },
    @{ Cell = "E45"; Value = "  -1.36%  "; IsPrice = $false },
    @{ Cell = "D46"; Value = "2.293"; IsPrice = $true },
    @{ Cell = "E46"; Value = "  +21.91%  "; IsPrice = $false },
    @{ Cell = "D47"; Value = "0.07175"; IsPrice = $true },
    @{ Cell = "E47"; Value = "  +4.60%  "; IsPrice = $false },
    @{ Cell = "B48"; Value = "EnergySwap"; IsPrice = $false },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; IsPrice = $false },
    @{ Cell = "D48"; Value = "11.74"; IsPrice = $true },
    @{ Cell = "E48"; Value = "  -3.40%  "; IsPrice = $false },
    @{ Cell = "B49"; Value = "Decentraland"; IsPrice = $false },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; IsPrice = $false },
    @{ Cell = "D49"; Value = "0.5343"; IsPrice = $true },
    @{ Cell = "E49"; Value = "  -4.88%  "; IsPrice = $false },
    @{ Cell = "D50"; Value = "1.861"; IsPrice = $true },
    @{ Cell = "E50"; Value = "  -5.98%  "; IsPrice = $false },
    @{ Cell = "D51"; Value = "109.20"; IsPrice = $true },
    @{ Cell = "E51"; Value = "  -2.26%  "; IsPrice = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.IsPrice) {
        # Force text format so numeric-looking price strings (e.g. "109.20", "1.002")
        # keep their exact text instead of being coerced into a Number and losing
        # trailing zeros / precision, matching the column's existing Text cell type.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
